# The deck originally ships two DrawingML themes:
#   ppt/theme/theme1.xml -> clrScheme "Office"      (name="Office Theme")
#   ppt/theme/theme2.xml -> clrScheme "Red Violet"   (name="Integral")
# theme2.xml is the one bound to the (single) slide master / the whole
# presentation's design, so it is what actually paints every slide.
# The target edit swaps the two themes' contents, i.e. the presentation's
# design should go from "Integral" (Red Violet) colors to the plain
# "Office Theme" colors. We reproduce that by rewriting the 12 theme
# color slots (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) of the
# active design's color scheme to the "Office Theme" RGB values.

function Hex-ToRGB([string]$hex) {
    # PowerPoint's RGB() / ColorFormat.RGB integers are 0x00BBGGRR
    # (red in the low byte), so build the value from R/G/B components.
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target "Office Theme" color scheme, in MsoThemeColorSchemeIndex order:
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink
$officeThemeHex = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$themeColors = $design.SlideMaster.Theme.ThemeColorScheme

for ($i = 1; $i -le $officeThemeHex.Count; $i++) {
    $themeColors.Item($i).RGB = Hex-ToRGB $officeThemeHex[$i - 1]
}
